$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "2020" year column, mirroring the existing "2019" column (P) formatting
$ws.Range("Q4").Value = 2020
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("Q5").Value = 90.6
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Active-cell selection saved with the file
[void]$ws.Range("P12").Select()
